$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0. Copy formatting of column AF into new column AG (header row 4 through row 61)
$ws.Range("AF4:AF61").Copy()
$ws.Range("AG4:AG61").PasteSpecial(-4122)  # xlPasteFormats

# 1. Add "Totals" header in AG4
$ws.Range("AG4").Value2 = "Totals"

# 2. Update header A4 "Crop years" -> "Country"
$ws.Range("A4").Value2 = "Country"

# 3. Update A6 "Bolivia (Plurinational State of)" -> "Bolivia"
$ws.Range("A6").Value2 = "Bolivia"

# 4. Row-total formulas: AG5 standalone, AG6:AG61 as one fill (shared formula)
$ws.Range("AG5").Formula = "=SUM(C5:AF5)"
$ws.Range("AG6:AG61").Formula = "=SUM(C6:AF6)"

# Row 60 is a blank spacer row - it must stay empty (no formula/value)
$ws.Range("AG60").ClearContents()

# AG61 carries style 15 (matches the blank spacer row above), not the bold-total row style
$ws.Range("AF60").Copy()
$ws.Range("AG61").PasteSpecial(-4122)
$ws.Range("AG61").Formula = "=SUM(C61:AF61)"

Write-Host "A4:" $ws.Range("A4").Value2
Write-Host "A6:" $ws.Range("A6").Value2
Write-Host "AG4:" $ws.Range("AG4").Value2
Write-Host "AG5:" $ws.Range("AG5").Value2
Write-Host "AG6:" $ws.Range("AG6").Value2
Write-Host "AG61:" $ws.Range("AG61").Value2
